$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 107162000000.0
$ws.Range("B4").Value = 4943000000.0
$ws.Range("G4").Value = 4097000000.0
$ws.Range("G5").Value = 12026000000.0
$ws.Range("G7").Value = 37031000000.0
$ws.Range("G8").Value = 99899000000.0
$ws.Range("G10").Value = 40457000000.0
$ws.Range("B13").Value = 16460000000.0
$ws.Range("G13").Value = 45111000000.0
$ws.Range("G15").Value = 15214000000.0
$ws.Range("B17").Value = 53255000000.0
$ws.Range("G17").Value = 36263000000.0
$ws.Range("G19").Value = 93078000000.0
$ws.Range("B21").Value = -149000000.0
$ws.Range("G22").Value = 55848000000.0
$ws.Range("B30").Value = 16686300000.0
$ws.Range("G32").Value = 1130000000.0
$ws.Range("G33").Value = 108292000000.0
